$d = $word.ActiveDocument

# 1) "...bis Landesliga" -> "...bis Hessenliga" (title + table cell, both occurrences)
$d.Content.Find.Execute("Landesliga", $true, $false, $false, $false, $false, $true, 1, $false, "Hessenliga", 2)

# 2) Address paragraph: re-set the "Move Sportwelt, Willy-Mock-Straße 2" span so the
#    stale spell-check proofErr markers around "Sportwelt" and "Straße" are cleared.
#    (A same-text assignment is a no-op, so stage through a placeholder first.)
$addr = $d.Content
$addr.Find.Execute("Move Sportwelt, Willy-Mock-Straße 2")
$addr.Text = "@@ADDR_PLACEHOLDER@@"
$addr2 = $d.Content
$addr2.Find.Execute("@@ADDR_PLACEHOLDER@@")
$addr2.Text = "Move Sportwelt, Willy-Mock-Straße 2"
